$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 200
$ws.Range('A200').Value = 198
$ws.Range('B200').Value = 6236254
$ws.Range('C200').Value = 'Venezuela Primera Division'
$ws.Range('D200').Value = 'Venezuela Primera Division'
$ws.Range('E200').Value = 45199.6875
$ws.Range('F200').Value = 'Academia Puerto Cabello'
$ws.Range('G200').Value = 'Estudiantes Merida'
$ws.Range('H200').Value = 1
$ws.Range('I200').Value = 0
$ws.Range('J200').Value = 'H'
$ws.Range('K200').Value = 1.727
$ws.Range('L200').Value = 3.4
$ws.Range('M200').Value = 4.333
$ws.Range('N200').Value = 1.666
$ws.Range('O200').Value = 3.4
$ws.Range('P200').Value = 4.75
$ws.Range('Q200').Value = -0.75
$ws.Range('R200').Value = 1.875
$ws.Range('S200').Value = 1.925
$ws.Range('T200').Value = 2.5
$ws.Range('U200').Value = 1.9
$ws.Range('V200').Value = 1.9
$ws.Range('W200').Value = 0.6659999999999999
$ws.Range('X200').Value = -1
$ws.Range('Y200').Value = -1
$ws.Range('Z200').Value = 0.4375
$ws.Range('AA200').Value = -0.5
$ws.Range('AB200').Value = -1
$ws.Range('AC200').Value = 0.8999999999999999

# Row 201
$ws.Range('A201').Value = 199
$ws.Range('B201').Value = 6236255
$ws.Range('C201').Value = 'Venezuela Primera Division'
$ws.Range('D201').Value = 'Venezuela Primera Division'
$ws.Range('E201').Value = 45199.6875
$ws.Range('F201').Value = 'Deportivo Rayo Zuliano'
$ws.Range('G201').Value = 'Caracas'
$ws.Range('H201').Value = 0
$ws.Range('I201').Value = 0
$ws.Range('J201').Value = 'D'
$ws.Range('K201').Value = 3.75
$ws.Range('L201').Value = 3.1
$ws.Range('M201').Value = 1.95
$ws.Range('N201').Value = 2.9
$ws.Range('O201').Value = 2.875
$ws.Range('P201').Value = 2.45
$ws.Range('Q201').Value = 0.25
$ws.Range('R201').Value = 1.775
$ws.Range('S201').Value = 2.025
$ws.Range('T201').Value = 2.25
$ws.Range('U201').Value = 1.85
$ws.Range('V201').Value = 1.95
$ws.Range('W201').Value = -1
$ws.Range('X201').Value = 1.875
$ws.Range('Y201').Value = -1
$ws.Range('Z201').Value = 0.3875
$ws.Range('AA201').Value = -0.5
$ws.Range('AB201').Value = -1
$ws.Range('AC201').Value = 0.95

# Row 202
$ws.Range('A202').Value = 200
$ws.Range('B202').Value = 6236251
$ws.Range('C202').Value = 'Venezuela Primera Division'
$ws.Range('D202').Value = 'Venezuela Primera Division'
$ws.Range('E202').Value = 45199.6875
$ws.Range('F202').Value = 'Angostura FC'
$ws.Range('G202').Value = 'Portuguesa'
$ws.Range('H202').Value = 1
$ws.Range('I202').Value = 2
$ws.Range('J202').Value = 'A'
$ws.Range('K202').Value = 3.1
$ws.Range('L202').Value = 3.2
$ws.Range('M202').Value = 2.15
$ws.Range('N202').Value = 4
$ws.Range('O202').Value = 3.6
$ws.Range('P202').Value = 1.75
$ws.Range('Q202').Value = 0.75
$ws.Range('R202').Value = 1.8
$ws.Range('S202').Value = 2
$ws.Range('T202').Value = 2.5
$ws.Range('U202').Value = 1.95
$ws.Range('V202').Value = 1.85
$ws.Range('W202').Value = -1
$ws.Range('X202').Value = -1
$ws.Range('Y202').Value = 0.75
$ws.Range('Z202').Value = -0.5
$ws.Range('AA202').Value = 0.5
$ws.Range('AB202').Value = 0.95
$ws.Range('AC202').Value = -1

# Row 203
$ws.Range('A203').Value = 201
$ws.Range('B203').Value = 6236253
$ws.Range('C203').Value = 'Venezuela Primera Division'
$ws.Range('D203').Value = 'Venezuela Primera Division'
$ws.Range('E203').Value = 45199.6875
$ws.Range('F203').Value = 'Deportivo La Guaira'
$ws.Range('G203').Value = 'UCV'
$ws.Range('H203').Value = 0
$ws.Range('I203').Value = 0
$ws.Range('J203').Value = 'D'
$ws.Range('K203').Value = 1.833
$ws.Range('L203').Value = 3.25
$ws.Range('M203').Value = 4
$ws.Range('N203').Value = 2
$ws.Range('O203').Value = 3.2
$ws.Range('P203').Value = 3.5
$ws.Range('Q203').Value = -0.25
$ws.Range('R203').Value = 1.775
$ws.Range('S203').Value = 2.025
$ws.Range('T203').Value = 2.25
$ws.Range('U203').Value = 1.9
$ws.Range('V203').Value = 1.9
$ws.Range('W203').Value = -1
$ws.Range('X203').Value = 2.2
$ws.Range('Y203').Value = -1
$ws.Range('Z203').Value = -0.5
$ws.Range('AA203').Value = 0.5125
$ws.Range('AB203').Value = -1
$ws.Range('AC203').Value = 0.8999999999999999

# Row 206
$ws.Range('A206').Value = 204
$ws.Range('B206').Value = 6236615
$ws.Range('C206').Value = 'Venezuela Primera Division'
$ws.Range('D206').Value = 'Venezuela Primera Division'
$ws.Range('E206').Value = 45206.6875
$ws.Range('F206').Value = 'Deportivo Rayo Zuliano'
$ws.Range('G206').Value = 'Academia Puerto Cabello'
$ws.Range('H206').Value = 1
$ws.Range('I206').Value = 0
$ws.Range('J206').Value = 'H'
$ws.Range('K206').Value = 2.375
$ws.Range('L206').Value = 3.3
$ws.Range('M206').Value = 2.625
$ws.Range('N206').Value = 2.45
$ws.Range('O206').Value = 3.2
$ws.Range('P206').Value = 2.55
$ws.Range('Q206').Value = 0
$ws.Range('R206').Value = 1.875
$ws.Range('S206').Value = 1.925
$ws.Range('T206').Value = 2.5
$ws.Range('U206').Value = 2
$ws.Range('V206').Value = 1.8
$ws.Range('W206').Value = 1.45
$ws.Range('X206').Value = -1
$ws.Range('Y206').Value = -1
$ws.Range('Z206').Value = 0.875
$ws.Range('AA206').Value = -1
$ws.Range('AB206').Value = -1
$ws.Range('AC206').Value = 0.8

# Row 207
$ws.Range('A207').Value = 205
$ws.Range('B207').Value = 6236616
$ws.Range('C207').Value = 'Venezuela Primera Division'
$ws.Range('D207').Value = 'Venezuela Primera Division'
$ws.Range('E207').Value = 45206.6875
$ws.Range('F207').Value = 'UCV'
$ws.Range('G207').Value = 'Metropolitanos FC'
$ws.Range('H207').Value = 3
$ws.Range('I207').Value = 2
$ws.Range('J207').Value = 'H'
$ws.Range('K207').Value = 3.3
$ws.Range('L207').Value = 3.2
$ws.Range('M207').Value = 2.05
$ws.Range('N207').Value = 2.75
$ws.Range('O207').Value = 3.2
$ws.Range('P207').Value = 2.3
$ws.Range('Q207').Value = 0.25
$ws.Range('R207').Value = 1.75
$ws.Range('S207').Value = 2.05
$ws.Range('T207').Value = 2.5
$ws.Range('U207').Value = 1.975
$ws.Range('V207').Value = 1.825
$ws.Range('W207').Value = 1.75
$ws.Range('X207').Value = -1
$ws.Range('Y207').Value = -1
$ws.Range('Z207').Value = 0.75
$ws.Range('AA207').Value = -1
$ws.Range('AB207').Value = 0.9750000000000001
$ws.Range('AC207').Value = -1

# Row 237
$ws.Range('A237').Value = 235
$ws.Range('B237').Value = 7842503
$ws.Range('C237').Value = 'Venezuela Primera Division'
$ws.Range('D237').Value = 'Venezuela Primera Division'
$ws.Range('E237').Value = 45339.66666666666
$ws.Range('F237').Value = 'Metropolitanos FC'
$ws.Range('G237').Value = 'Zamora'
$ws.Range('K237').Value = 1.8
$ws.Range('L237').Value = 3.5
$ws.Range('M237').Value = 3.75
$ws.Range('N237').Value = 1.85
$ws.Range('O237').Value = 3.5
$ws.Range('P237').Value = 3.6
$ws.Range('Q237').Value = -0.5
$ws.Range('R237').Value = 1.9
$ws.Range('S237').Value = 1.9
$ws.Range('T237').Value = 2.5
$ws.Range('U237').Value = 2.025
$ws.Range('V237').Value = 1.775
$ws.Range('W237').Value = 0
$ws.Range('X237').Value = 0
$ws.Range('Y237').Value = 0
$ws.Range('Z237').Value = 0
$ws.Range('AA237').Value = 0

# Row 238
$ws.Range('A238').Value = 236
$ws.Range('B238').Value = 7842504
$ws.Range('C238').Value = 'Venezuela Primera Division'
$ws.Range('D238').Value = 'Venezuela Primera Division'
$ws.Range('E238').Value = 45339.78125
$ws.Range('F238').Value = 'Angostura FC'
$ws.Range('G238').Value = 'Deportivo La Guaira'
$ws.Range('K238').Value = 2.75
$ws.Range('L238').Value = 3
$ws.Range('M238').Value = 2.45
$ws.Range('N238').Value = 3.1
$ws.Range('O238').Value = 3
$ws.Range('P238').Value = 2.15
$ws.Range('Q238').Value = 0.25
$ws.Range('R238').Value = 1.925
$ws.Range('S238').Value = 1.875
$ws.Range('T238').Value = 2.25
$ws.Range('U238').Value = 1.975
$ws.Range('V238').Value = 1.825
$ws.Range('W238').Value = 0
$ws.Range('X238').Value = 0
$ws.Range('Y238').Value = 0
$ws.Range('Z238').Value = 0
$ws.Range('AA238').Value = 0

# Row 239
$ws.Range('A239').Value = 237
$ws.Range('B239').Value = 7842507
$ws.Range('C239').Value = 'Venezuela Primera Division'
$ws.Range('D239').Value = 'Venezuela Primera Division'
$ws.Range('E239').Value = 45339.78125
$ws.Range('F239').Value = 'Academia Puerto Cabello'
$ws.Range('G239').Value = 'Estudiantes Merida'
$ws.Range('K239').Value = 1.727
$ws.Range('L239').Value = 3.5
$ws.Range('M239').Value = 4.2
$ws.Range('N239').Value = 1.75
$ws.Range('O239').Value = 3.6
$ws.Range('P239').Value = 4
$ws.Range('Q239').Value = -0.75
$ws.Range('R239').Value = 1.975
$ws.Range('S239').Value = 1.825
$ws.Range('T239').Value = 2.5
$ws.Range('U239').Value = 1.9
$ws.Range('V239').Value = 1.9
$ws.Range('W239').Value = 0
$ws.Range('X239').Value = 0
$ws.Range('Y239').Value = 0
$ws.Range('Z239').Value = 0
$ws.Range('AA239').Value = 0

# Row 240
$ws.Range('A240').Value = 238
$ws.Range('B240').Value = 7842505
$ws.Range('C240').Value = 'Venezuela Primera Division'
$ws.Range('D240').Value = 'Venezuela Primera Division'
$ws.Range('E240').Value = 45339.89583333334
$ws.Range('F240').Value = 'Deportivo Tachira'
$ws.Range('G240').Value = 'Deportivo Rayo Zuliano'
$ws.Range('K240').Value = 1.444
$ws.Range('L240').Value = 3.75
$ws.Range('M240').Value = 7
$ws.Range('N240').Value = 1.5
$ws.Range('O240').Value = 3.6
$ws.Range('P240').Value = 6
$ws.Range('Q240').Value = -1
$ws.Range('R240').Value = 1.875
$ws.Range('S240').Value = 1.925
$ws.Range('T240').Value = 2.5
$ws.Range('U240').Value = 1.85
$ws.Range('V240').Value = 1.95
$ws.Range('W240').Value = 0
$ws.Range('X240').Value = 0
$ws.Range('Y240').Value = 0
$ws.Range('Z240').Value = 0
$ws.Range('AA240').Value = 0

# Row 241
$ws.Range('A241').Value = 239
$ws.Range('B241').Value = 7842502
$ws.Range('C241').Value = 'Venezuela Primera Division'
$ws.Range('D241').Value = 'Venezuela Primera Division'
$ws.Range('E241').Value = 45340.66666666666
$ws.Range('F241').Value = 'CD Hermanos Colmenares'
$ws.Range('G241').Value = 'UCV'
$ws.Range('K241').Value = 2.3
$ws.Range('L241').Value = 3.1
$ws.Range('M241').Value = 2.875
$ws.Range('N241').Value = 2.4
$ws.Range('O241').Value = 3.1
$ws.Range('P241').Value = 2.75
$ws.Range('Q241').Value = -0.25
$ws.Range('R241').Value = 2.1
$ws.Range('S241').Value = 1.7
$ws.Range('T241').Value = 2
$ws.Range('U241').Value = 1.875
$ws.Range('V241').Value = 1.925
$ws.Range('W241').Value = 0
$ws.Range('X241').Value = 0
$ws.Range('Y241').Value = 0
$ws.Range('Z241').Value = 0
$ws.Range('AA241').Value = 0

# Row 242
$ws.Range('A242').Value = 240
$ws.Range('B242').Value = 7842506
$ws.Range('C242').Value = 'Venezuela Primera Division'
$ws.Range('D242').Value = 'Venezuela Primera Division'
$ws.Range('E242').Value = 45340.78125
$ws.Range('F242').Value = 'Monagas'
$ws.Range('G242').Value = 'Carabobo'
$ws.Range('K242').Value = 2.8
$ws.Range('L242').Value = 3
$ws.Range('M242').Value = 2.4
$ws.Range('N242').Value = 2.8
$ws.Range('O242').Value = 3
$ws.Range('P242').Value = 2.4
$ws.Range('Q242').Value = 0
$ws.Range('R242').Value = 2.05
$ws.Range('S242').Value = 1.75
$ws.Range('T242').Value = 2.25
$ws.Range('U242').Value = 1.975
$ws.Range('V242').Value = 1.825
$ws.Range('W242').Value = 0
$ws.Range('X242').Value = 0
$ws.Range('Y242').Value = 0
$ws.Range('Z242').Value = 0
$ws.Range('AA242').Value = 0

# Delete row 243 (removed in this update)
$ws.Rows(243).Delete()

Write-Host "Update complete"